# Generate Report for handoff
#
# - Row that used to describe 76f98c9c-...md / "Handoff transform failed"
#   is now 3bac9dce-...md / "Ready for handoff" (status moved on).
# - A brand-new source file ffff75bd9d4b-...md shows up with its own
#   "Ready for handoff" row, pushing .localization-config down a row.
# - The two per-language sheets pick up the freshly generated handoff
#   package (.xlf) with its timestamp + "Include" dependency reason.

$wb = $excel.ActiveWorkbook

$oldMdName = "76f98c9c-c65b-44b5-9710-80dc295ad7b5.md"
$newMdName = "3bac9dce-e042-4ff9-be51-261a1d015ce5.md"
$newMdName2 = "ffff75bd9d4b-df6f-4eb7-b85d-5406951d8eb4.md"
$cfgName = ".localization-config"

$zhXlf = "3bac9dce-e042-4ff9-be51-261a1d015ce5.0343cc1b35dc03dcaf1978a748210400f6a412f7.zh-cn.xlf"
$deXlf = "3bac9dce-e042-4ff9-be51-261a1d015ce5.0343cc1b35dc03dcaf1978a748210400f6a412f7.de-de.xlf"

$zhHandoffTime = "2016-02-15 08:44:57"
$deHandoffTime = "2016-02-15 08:45:12"
$epoch = "0001-01-01 00:00:00"

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8382ec6132d6593ed9f911cbe761e03a9deae810"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4047eb00623a2b5c725bd6799689b239b1fa98bb/.localization-config"
$xlfBaseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8382ec6132d6593ed9f911cbe761e03a9deae810/e2e"

function Update-OverviewSheet($ws) {

    # Drop the old hyperlinks so Excel renumbers the relationship ids
    # cleanly from scratch, in the same order they appear on the sheet.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("A3").Hyperlinks.Delete()

    $ws.Range("A2").Value = $newMdName
    $ws.Range("B2").Value = "Ready for handoff"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newMdName", "", "", $newMdName)

    $ws.Range("A3").Value = $newMdName2
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/e2e/$newMdName2", "", "", $newMdName2)

    $ws.Range("A4").Value = $cfgName
    $ws.Range("B4").Value = "Not to be localized"
    $ws.Range("C4").Value = "Not to be localized"
    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)
}

function Update-LanguageSheet($ws, $xlfName, $handoffTime) {

    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("A3").Hyperlinks.Delete()

    # Row 2: the renamed/progressed source file
    $ws.Range("A2").Value = $newMdName
    $ws.Range("B2").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newMdName", "", "", $newMdName)

    $ws.Range("C2").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("C2"), "$xlfBaseUrl/$xlfName", "", "", $xlfName)

    $ws.Range("D2").Value = $handoffTime
    $ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = "Include"

    # Row 3: brand-new source file, same handoff package dependency
    $ws.Range("A3").Value = $newMdName2
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/e2e/$newMdName2", "", "", $newMdName2)

    $ws.Range("C3").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("C3"), "$xlfBaseUrl/$xlfName", "", "", $xlfName)

    $ws.Range("D3").Value = $handoffTime
    $ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = "Include"

    # Row 4: .localization-config, pushed down from row 3
    $ws.Range("A4").Value = $cfgName
    $ws.Range("B4").Value = "Not to be localized"
    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)

    $ws.Range("D4").Value = $epoch
    $ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ws.Range("G4").Value = $epoch
    $ws.Range("H4").Value = "Ignored"
}

Update-OverviewSheet $wb.Worksheets.Item("Overview")
Update-LanguageSheet $wb.Worksheets.Item("zh-cn") $zhXlf $zhHandoffTime
Update-LanguageSheet $wb.Worksheets.Item("de-de") $deXlf $deHandoffTime
